$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 280, shifting rows 280:378 down to 281:379
$ws.Rows.Item(280).Insert()

# Populate the newly inserted row 280 with the new record's data
$ws.Cells.Item(280, 1).Value = 7
$ws.Cells.Item(280, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(280, 3).Value = "Ñuble"
$ws.Cells.Item(280, 4).Value = 44588
$ws.Cells.Item(280, 5).Value = 16
$ws.Cells.Item(280, 6).Value = 100112020
$ws.Cells.Item(280, 7).Value = "Tomate"
$ws.Cells.Item(280, 8).Value = "Larga vida"
$ws.Cells.Item(280, 9).Value = "Primera"
$ws.Cells.Item(280, 10).Value = 800
$ws.Cells.Item(280, 11).Value = 5500
$ws.Cells.Item(280, 12).Value = 6000
$ws.Cells.Item(280, 13).Value = 5750
$ws.Cells.Item(280, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(280, 15).Value = "Región del Maule"
$ws.Cells.Item(280, 16).Value = 383
$ws.Cells.Item(280, 17).Value = 15
$ws.Cells.Item(280, 18).Value = "Hortaliza"

# Match the date-column style (s="2") used by the rest of column D
$ws.Cells.Item(280, 4).NumberFormat = $ws.Cells.Item(281, 4).NumberFormat
